$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.317.43"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "3.680.90"
$ws.Range("E3").Value = "  -3.43%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "683.05"
$ws.Range("E5").Value = "  -3.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.53"
$ws.Range("E6").Value = "  -4.68%  "
$ws.Range("D7").Value = "3.680.56"
$ws.Range("E7").Value = "  -3.45%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  -4.30%  "
$ws.Range("E10").Value = "  -7.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.24"
$ws.Range("E11").Value = "  -4.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.451"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("E13").Value = "  -6.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.57"
$ws.Range("E14").Value = "  -6.70%  "
$ws.Range("D15").Value = "4.303.90"
$ws.Range("E15").Value = "  -3.37%  "
$ws.Range("D16").Value = "3.685.40"
$ws.Range("E16").Value = "  -3.76%  "
$ws.Range("D17").Value = "69.351.70"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.32"
$ws.Range("E19").Value = "  -5.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.66"
$ws.Range("E20").Value = "  -6.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "479.58"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("E22").Value = "  -7.78%  "
$ws.Range("E23").Value = "  -8.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.06"
$ws.Range("E24").Value = "  -4.97%  "
$ws.Range("D25").Value = "3.828.30"
$ws.Range("E25").Value = "  -3.38%  "
$ws.Range("E26").Value = "  -11.06%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.50"
$ws.Range("E28").Value = "  -4.87%  "
$ws.Range("E29").Value = "  -7.65%  "
$ws.Range("E30").Value = "  -10.49%  "
$ws.Range("E31").Value = "  -10.79%  "
$ws.Range("E32").Value = "  -5.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.89"
$ws.Range("E33").Value = "  -6.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.13"
$ws.Range("E34").Value = "  -6.76%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.163"
$ws.Range("E36").Value = "  -4.77%  "
$ws.Range("D37").Value = "3.644.86"
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("E38").Value = "  -5.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.14"
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0942"
$ws.Range("E40").Value = "  -7.23%  "
$ws.Range("E42").Value = "  -6.19%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("E44").Value = "  -7.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.22"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "154.35"
$ws.Range("E46").Value = "  -7.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.84"
$ws.Range("E47").Value = "  -12.17%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000282"
$ws.Range("E49").Value = "  -12.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "392.14"
$ws.Range("E50").Value = "  -7.50%  "
$ws.Range("E51").Value = "  -6.00%  "
